# Rename the ambiguous "Cross references" sheet to "Database references"
# and make it the active sheet (mirrors the tabSelected/activeTab move in
# the target diff), matching the commit message:
#   "in all models change Cross references to Database references to
#    replace an ambiguous name with a clear one"

$wb = $excel.ActiveWorkbook

# Rename the sheet (sheetId="12", previously "Cross references").
$ws = $wb.Worksheets.Item("Cross references")
$ws.Name = "Database references"

# Make it the active/selected sheet - this moves `tabSelected="1"` off the
# previously-selected sheet ("Taxon") and onto this one, and updates the
# workbook's bookViews/activeTab accordingly.
$ws.Activate()

# Enable iterative calculation with a tighter convergence delta
# (<calcPr iterateDelta="1E-4" .../> in the target workbook).
$excel.Iteration = $true
$excel.MaxChange = 0.0001
